$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F (想去人数 / want-to-go count)
$updates = @{
    2  = 245
    3  = 272
    4  = 282
    6  = 278
    7  = 6691
    11 = 81
    15 = 19
    17 = 563
    18 = 62
}

# Both "展览" and "全部类型" sheets contain the same data and need the same update
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F" + $row).Value = $updates[$row]
    }
}
